$d = $word.ActiveDocument

function Split-TextIntoRuns($findText, $pieceLengths) {
    # Locate the target text first so we work with absolute character
    # offsets into the document story.
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    $base = $rng.Start

    # Toggling (and reverting) a character formatting property on each
    # non-empty sub-range forces Word to break run boundaries at the
    # edges of that sub-range, without changing any visible formatting,
    # so the run ends up split exactly where the pieces join.
    $pos = $base
    foreach ($len in $pieceLengths) {
        $piece = $d.Range($pos, $pos + $len)
        $piece.Font.Bold = 1
        $piece.Font.Bold = 0
        $pos = $pos + $len
    }
}

# "#Career:1#" -> "#Career" + ":1" + "#"
Split-TextIntoRuns "#Career:1#" @(7, 2, 1)

# "#Education:1#" -> "#Education:" + "1" + "#"
Split-TextIntoRuns "#Education:1#" @(11, 1, 1)
